$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new values are unambiguous text (safe to set directly) ---
$ws.Range("D2").Value = "62.244.39"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").Value = "3.002.33"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("E6").Value = "  -5.87%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "2.998.43"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("E10").Value = "  -5.41%  "
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("E14").Value = "  -6.51%  "
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "3.493.51"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "62.227.29"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").Value = "3.001.19"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  -5.03%  "
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E25").Value = "  -4.90%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E26").Value = "  -10.49%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -7.02%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("E31").Value = "  -6.29%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E33").Value = "  -7.21%  "
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0794"
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("E37").Value = "  -5.18%  "
$ws.Range("E38").Value = "  -5.98%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("E41").Value = "  -11.60%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("E43").Value = "  -10.53%  "
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("E45").Value = "  -7.09%  "
$ws.Range("D46").Value = "2.721.25"
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("E47").Value = "  -5.78%  "
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E51").Value = "  -2.47%  "

# --- Cells whose new values look numeric; force text storage to preserve exact formatting ---
$textCells = @("D4", "D5", "D6", "D8", "D12", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D39", "D40", "D41", "D43", "D45", "D47", "D48", "D51")
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "582.31"
$ws.Range("D6").Value = "145.68"
$ws.Range("D8").Value = "0.526"
$ws.Range("D12").Value = "0.460"
$ws.Range("D14").Value = "34.41"
$ws.Range("D17").Value = "7.13"
$ws.Range("D20").Value = "459.85"
$ws.Range("D21").Value = "13.97"
$ws.Range("D22").Value = "0.686"
$ws.Range("D23").Value = "7.44"
$ws.Range("D24").Value = "81.47"
$ws.Range("D25").Value = "12.28"
$ws.Range("D26").Value = "2.21"
$ws.Range("D28").Value = "9.95"
$ws.Range("D29").Value = "0.998"
$ws.Range("D30").Value = "2.62"
$ws.Range("D31").Value = "7.02"
$ws.Range("D32").Value = "28.32"
$ws.Range("D33").Value = "2.08"
$ws.Range("D35").Value = "1.03"
$ws.Range("D37").Value = "5.76"
$ws.Range("D39").Value = "9.24"
$ws.Range("D40").Value = "50.25"
$ws.Range("D41").Value = "2.86"
$ws.Range("D43").Value = "393.57"
$ws.Range("D45").Value = "0.271"
$ws.Range("D47").Value = "36.60"
$ws.Range("D48").Value = "128.30"
$ws.Range("D51").Value = "2.20"
foreach ($cellref in $textCells) {
    $ws.Range($cellref).ClearFormats()
}
